$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "Child " and "argument used to group together individual nodes in
# complex systems (e.g., Henry's reach)" were two separate runs; collapse
# them into a single run (text content is unchanged, only the run split
# is simplified). A Find/Replace over the full phrase causes Word to
# rewrite the matched span as one run.
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$find1.Execute("Child argument used to group together individual nodes in complex systems (e.g., Henry" + [char]8217 + "s reach)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Child argument used to group together individual nodes in complex systems (e.g., Henry" + [char]8217 + "s reach)", `
    2) | Out-Null

# --- Change 2 -----------------------------------------------------------
# Append a new bullet item (same list - numId 27, level 0) after the
# final "TagObs_Directionality_YY-MM-DD.csv" bullet, describing the new
# "TagObs_FinalPaths_YY-MM-DD.csv" output file.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.ListFormat.ListLevelNumber = 1

$dq = [char]8220
$rq = [char]8221

$newPara.Range.Text = "There is a code chunk following the primary workflow that subsets the " + $dq + "TagObs_Directionality" + $rq + " dataset to the final detections for each tag and writes it to " + $dq + "TagObs_FinalPaths_YY-MM-DD.csv" + $rq + ". This provides the complete detection path for each fish through nodes specified in the " + $dq + "node_direction" + $rq + " file."
